$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 30 ("N" row), pushing it down to row 32
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# Helper data for the two new rows (row number -> ordered values for columns A..L)
$row30 = @("Cognitive Difficulty","0.12","0.16","0.14","0.15","0.1","0.08","0.1","0.12","0.14","0.09","0.12")
$row31 = @("Independence Difficulty","0.18","0.22","0.19","0.19","0.15","0.13","0.16","0.18","0.21","0.14","0.17")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "30")
    $cell.NumberFormat = "@"
    $cell.Value = $row30[$i]
    $cell.Style = "Normal"
}

for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "31")
    $cell.NumberFormat = "@"
    $cell.Value = $row31[$i]
    $cell.Style = "Normal"
}
